$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old F1 header ("ElementName3") so the used range shrinks from A1:F13 to A1:E13
$ws.Range("F1").ClearContents()

# Fill in the new B:E columns (Name / ChefModule / ElementName1 / ElementName2) for rows 2-13.
# Column A (Code, GSTR11..GSTR26) already holds the right values and is left untouched.
$ws.Range("B2").Value = "pede. Suspendisse dui."
$ws.Range("C2").Value = "EL Haddad"
$ws.Range("D2").Value = "Nullam feugiat placerat"
$ws.Range("E2").Value = "varius et, euismod"

$ws.Range("B3").Value = "a nunc. In"
$ws.Range("C3").Value = "Badir"
$ws.Range("D3").Value = "sodales nisi magna"
$ws.Range("E3").Value = "elementum sem, vitae"

$ws.Range("B4").Value = "amet metus. Aliquam"
$ws.Range("C4").Value = "Ezzine"
$ws.Range("D4").Value = "Cras vulputate velit"
$ws.Range("E4").Value = "scelerisque neque sed"

$ws.Range("B5").Value = "quam vel sapien"
$ws.Range("C5").Value = "El Alami Hassoun"
$ws.Range("D5").Value = "Nunc mauris elit,"
$ws.Range("E5").Value = "libero et tristique"

$ws.Range("B6").Value = "feugiat nec, diam."
$ws.Range("C6").Value = "Lazaar"
$ws.Range("D6").Value = "pellentesque. Sed dictum."
$ws.Range("E6").Value = "ridiculus mus. Proin"

$ws.Range("B7").Value = "nonummy. Fusce fermentum"
$ws.Range("C7").Value = "El Haddad"
$ws.Range("D7").Value = "neque pellentesque massa"
$ws.Range("E7").Value = "Mauris eu turpis."

$ws.Range("B8").Value = "a, arcu. Sed"
$ws.Range("C8").Value = "EL Haddad"
$ws.Range("D8").Value = "sit amet risus."
$ws.Range("E8").Value = "Nulla facilisi. Sed"

$ws.Range("B9").Value = "Suspendisse eleifend. Cras"
$ws.Range("C9").Value = "El Alami Hassoun"
$ws.Range("D9").Value = "velit dui, semper"
$ws.Range("E9").Value = "ligula elit, pretium"

$ws.Range("B10").Value = "ante. Nunc mauris"
$ws.Range("C10").Value = "Badir"
$ws.Range("D10").Value = "tortor at risus."
$ws.Range("E10").Value = "felis. Donec tempor,"

$ws.Range("B11").Value = "lobortis quam a"
$ws.Range("C11").Value = "Ezzine"
$ws.Range("D11").Value = "euismod est arcu"
$ws.Range("E11").Value = "ligula eu enim."

$ws.Range("B12").Value = "rhoncus. Nullam velit"
$ws.Range("C12").Value = "Ben Achrab"
$ws.Range("D12").Value = "ut dolor dapibus"
$ws.Range("E12").Value = "commodo tincidunt nibh."

$ws.Range("B13").Value = "Donec tincidunt. Donec"
$ws.Range("C13").Value = "EL Haddad"
$ws.Range("D13").Value = "ornare tortor at"
$ws.Range("E13").Value = "ac, feugiat non,"

# B2 carries an explicit black font color in the target workbook (new font/style entry)
$ws.Range("B2").Font.Color = 0

# Resize the data columns to fit the new (wider) content (values chosen so the
# persisted column width, after Excel's char->pixel rounding, lands as close as
# possible to the target widths of 27 / 16.140625 / 25.140625 / 23.5703125)
$ws.Columns.Item(2).ColumnWidth = 26.17
$ws.Columns.Item(3).ColumnWidth = 15.33
$ws.Columns.Item(4).ColumnWidth = 24.33
$ws.Columns.Item(5).ColumnWidth = 22.67

Write-Host "done"
